$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- Row 2: shift manager's own checkout (tair/hadad) gets corrected checkout time ---
Set-TextCell 2 8 "Fri, 04 Jan 2019 11:07:04"

# --- Row 4: yoni/machluf checkout corrected; total seconds recalculated ---
Set-TextCell 4 8 "Fri, 04 Jan 2019 10:44:37"
$ws.Cells.Item(4, 9).Value = -1419

# --- Row 5: tair/hadad (second shift) checkout corrected ---
Set-TextCell 5 8 "Fri, 04 Jan 2019 11:07:04"

# --- Row 10: replace the old open/incomplete "michal/tsho" presence record
#     with a new completed "tair/hadad" presence record ---
Set-TextCell 10 2 "tair"
Set-TextCell 10 3 "hadad"
Set-TextCell 10 4 "Fri, 04 Jan 2019 11:03:40"
Set-TextCell 10 5 "1"
Set-TextCell 10 6 "4"
Set-TextCell 10 7 "6"
Set-TextCell 10 8 "Fri, 04 Jan 2019 11:07:04"
$ws.Cells.Item(10, 9).Value = 204

# --- Row 11 (new): "michal/tsho" sign-in, still present (no checkout yet) ---
$ws.Cells.Item(11, 1).Value = 10
Set-TextCell 11 2 "michal"
Set-TextCell 11 3 "tsho"
Set-TextCell 11 4 "Fri, 04 Jan 2019 11:07:17"
Set-TextCell 11 5 "1"
Set-TextCell 11 6 "4"
Set-TextCell 11 7 "6"

$wb.Save()
